$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("optimization_parameters")

# Make "optimization_parameters" the active sheet (moves tabSelected / activeTab)
$ws.Activate()

# Insert a new row above row 9 (current row 9 "estimate_params" shifts down to row 10, etc.)
$ws.Rows.Item(9).Insert()

# Row 8 now represents the renamed "Model" -> "production_function" parameter
$ws.Range("A8").Value = "production_function"

# New row 9 holds the "L_curve" parameter with a default value of 0
$ws.Range("A9").Value = "L_curve"
$ws.Range("B9").Value = 0

# Update the visible selection to the new row
$ws.Range("A9:B9").Select()
